$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$value) {
    # Force the cell to remain text-typed (matches source data where numeric-
    # looking strings like "3.11" or multi-dot strings like "56.602.36" are
    # stored as text), then restore the original style so no stray formatting
    # is introduced.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextCell $ws.Range("D2") "56.602.36"
$ws.Range("E2").Value = "  +10.93%  "
Set-TextCell $ws.Range("D3") "3.250.99"
$ws.Range("E3").Value = "  +6.69%  "
$ws.Range("E4").Value = "  +0.13%  "
Set-TextCell $ws.Range("D5") "398.97"
$ws.Range("E5").Value = "  +3.31%  "
Set-TextCell $ws.Range("D6") "111.41"
$ws.Range("E6").Value = "  +9.35%  "
$ws.Range("E7").Value = "  +5.47%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +8.01%  "
Set-TextCell $ws.Range("D10") "39.49"
$ws.Range("E10").Value = "  +8.24%  "
Set-TextCell $ws.Range("D11") "0.0952"
$ws.Range("E11").Value = "  +12.49%  "
$ws.Range("E12").Value = "  +2.67%  "
Set-TextCell $ws.Range("D13") "3.759.23"
$ws.Range("E13").Value = "  +6.49%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws.Range("D14") "19.21"
$ws.Range("E14").Value = "  +5.33%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws.Range("D15") "8.12"
$ws.Range("E15").Value = "  +6.14%  "
Set-TextCell $ws.Range("D16") "3.249.15"
$ws.Range("E16").Value = "  +6.51%  "
$ws.Range("E17").Value = "  +7.17%  "
Set-TextCell $ws.Range("D18") "11.05"
$ws.Range("E18").Value = "  +3.86%  "
Set-TextCell $ws.Range("D19") "56.533.72"
$ws.Range("E19").Value = "  +10.84%  "
Set-TextCell $ws.Range("D20") "3.33"
$ws.Range("E20").Value = "  +4.45%  "
Set-TextCell $ws.Range("D21") "0.0000105"
$ws.Range("E21").Value = "  +9.86%  "
Set-TextCell $ws.Range("D22") "13.10"
$ws.Range("E22").Value = "  +7.49%  "
Set-TextCell $ws.Range("D23") "299.88"
$ws.Range("E23").Value = "  +13.78%  "
Set-TextCell $ws.Range("D24") "75.66"
$ws.Range("E24").Value = "  +9.01%  "
$ws.Range("E25").Value = "  +4.44%  "
Set-TextCell $ws.Range("D26") "8.19"
$ws.Range("E26").Value = "  +4.19%  "
Set-TextCell $ws.Range("D27") "28.28"
$ws.Range("E27").Value = "  +5.25%  "
$ws.Range("E28").Value = "  +4.32%  "
$ws.Range("E29").Value = "  +2.98%  "
$ws.Range("E30").Value = "  +5.92%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +6.77%  "
Set-TextCell $ws.Range("D33") "11.13"
$ws.Range("E33").Value = "  +7.64%  "
Set-TextCell $ws.Range("D34") "37.05"
$ws.Range("E34").Value = "  +5.51%  "
Set-TextCell $ws.Range("D35") "0.0489"
$ws.Range("E35").Value = "  +4.74%  "
$ws.Range("E36").Value = "  +3.10%  "
$ws.Range("E37").Value = "  +3.05%  "
$ws.Range("E38").Value = "  +5.36%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws.Range("D39") "3.11"
$ws.Range("E39").Value = "  +26.73%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws.Range("D40") "0.999"
$ws.Range("E40").Value = "  +0.00%  "
Set-TextCell $ws.Range("D41") "17.69"
$ws.Range("E41").Value = "  +8.58%  "
Set-TextCell $ws.Range("D42") "134.82"
$ws.Range("E42").Value = "  +4.18%  "
Set-TextCell $ws.Range("D43") "1.93"
$ws.Range("E43").Value = "  +6.16%  "
Set-TextCell $ws.Range("D44") "4.00"
$ws.Range("E44").Value = "  +5.97%  "
$ws.Range("E45").Value = "  +5.06%  "
$ws.Range("E46").Value = "  -1.18%  "
Set-TextCell $ws.Range("D47") "22.31"
$ws.Range("E47").Value = "  +4.00%  "
Set-TextCell $ws.Range("D48") "2.20"
$ws.Range("E48").Value = "  +58.78%  "
Set-TextCell $ws.Range("D49") "2.144.97"
$ws.Range("E49").Value = "  +4.67%  "
Set-TextCell $ws.Range("D50") "2.09"
$ws.Range("E50").Value = "  +0.46%  "
Set-TextCell $ws.Range("D51") "2.41"
$ws.Range("E51").Value = "  -2.27%  "
